$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44391, 'Primera', 400, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44742, 'Primera', 300, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44238, 'Primera', 300, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44238, 'Segunda', 200, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Provincia de Limarí', 556),
    @(44238, 'Tercera', 50, 8000, 8000, 8000, '$/bandeja 18 kilos', 'Provincia de Limarí', 444),
    @(44631, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44748, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44229, 'Primera', 200, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44627, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44726, 'Primera', 300, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44719, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44628, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44249, 'Primera', 400, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44249, 'Segunda', 200, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Provincia de Limarí', 556),
    @(44746, 'Primera', 500, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44784, 'Primera', 500, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44596, 'Primera', 150, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44383, 'Primera', 300, 16000, 16000, 16000, '$/bandeja 18 kilos', 'Provincia de Limarí', 889),
    @(44383, 'Segunda', 200, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44804, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44235, 'Primera', 400, 13000, 13000, 13000, '$/bandeja 18 kilos', 'Provincia de Limarí', 722),
    @(44235, 'Segunda', 200, 11000, 11000, 11000, '$/bandeja 18 kilos', 'Provincia de Limarí', 611),
    @(44235, 'Tercera', 100, 9000, 9000, 9000, '$/bandeja 18 kilos', 'Provincia de Limarí', 500),
    @(44803, 'Primera', 350, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44243, 'Especial', 300, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44243, 'Primera', 300, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Provincia de Limarí', 556),
    @(44243, 'Segunda', 150, 8000, 8000, 8000, '$/bandeja 18 kilos', 'Provincia de Limarí', 444),
    @(44721, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44614, 'Primera', 300, 15000, 15000, 15000, '$/caja 18 kilos granel', 'Provincia de Limarí', 833),
    @(44753, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44754, 'Primera', 400, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44635, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44802, 'Primera', 500, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44585, 'Primera', 200, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44750, 'Primera', 200, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44783, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44396, 'Primera', 250, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44396, 'Segunda', 150, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44777, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44729, 'Primera', 300, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44791, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44740, 'Primera', 400, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44725, 'Primera', 400, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Región de Arica y Parinacota', 833),
    @(44630, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44762, 'Primera', 300, 16000, 16000, 16000, '$/bandeja 18 kilos', 'Provincia de Limarí', 889),
    @(44245, 'Primera', 300, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Provincia de Limarí', 667),
    @(44245, 'Segunda', 200, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Provincia de Limarí', 556),
    @(44757, 'Primera', 300, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833),
    @(44806, 'Primera', 200, 14000, 14000, 14000, '$/bandeja 18 kilos', 'Provincia de Limarí', 778),
    @(44736, 'Primera', 200, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí', 833)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 4).Value = $row[0]   # D
    $ws.Cells.Item($r, 9).Value = $row[1]   # I
    $ws.Cells.Item($r, 10).Value = $row[2]  # J
    $ws.Cells.Item($r, 11).Value = $row[3]  # K
    $ws.Cells.Item($r, 12).Value = $row[4]  # L
    $ws.Cells.Item($r, 13).Value = $row[5]  # M
    $ws.Cells.Item($r, 14).Value = $row[6]  # N
    $ws.Cells.Item($r, 15).Value = $row[7]  # O
    $ws.Cells.Item($r, 16).Value = $row[8]  # P
}
